# Ajout de l'activite de la journee (nouvelle ligne 65 dans le journal de travail)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing entry (row 64) onto the new
# row 65, cell by cell, so the new cells pick up the same styles
# (date format for column A, wrap-text for column B) already present
# in the shared style table instead of creating new ones.
$ws.Cells.Item(64, 1).Copy() | Out-Null
$ws.Cells.Item(65, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Cells.Item(64, 2).Copy() | Out-Null
$ws.Cells.Item(65, 2).PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Cells.Item(64, 3).Copy() | Out-Null
$ws.Cells.Item(65, 3).PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Fill in the new journal entry.
$ws.Cells.Item(65, 1).Value = 43209

$newText = "Finalisation de la doc. J'ai ajouté des captures d'écrans de l'application mobile pour avoir une preuve qu'elle fonctionne. J'ai ajouté des trucs dans certains chapitres et fais attention à l'horographe, même si je pense avoir laissé des fautes`n"
$ws.Cells.Item(65, 2).Value = $newText
$ws.Cells.Item(65, 2).WrapText = $true

$ws.Cells.Item(65, 3).Value = "2 périodes"

$ws.Rows.Item(65).RowHeight = 60

# Reflect the post-entry selection (user landed on B66 after typing the row).
$ws.Range("B66").Select()
